$d = $word.ActiveDocument

# The document currently ends with a single empty trailing paragraph
# immediately before the sectPr. We need to insert a whole new block of
# meeting-minute content (one Title paragraph, one Heading1 "Attendees"
# paragraph, and eight ListParagraph/numId=2 bullet paragraphs) right
# before that trailing empty paragraph, followed by one brand-new empty
# paragraph - leaving the original trailing empty paragraph untouched,
# immediately before the sectPr, exactly as in the target revision.

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# A collapsed (zero-length) range at the very start of that last, empty
# paragraph. Inserting raw WordprocessingML at a collapsed range splices
# in whole new paragraphs for everything up to the final fragment
# paragraph; since our fragment's final paragraph is itself empty, it
# becomes a fresh empty paragraph of its own rather than merging text
# into - or consuming - the pre-existing trailing empty paragraph.
$insertRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$newContentXml = @'
<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>Thursday 31/03</w:t></w:r><w:r><w:t xml:space="preserve"> 02:20pm</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Attendees: Mab, Andre</w:t></w:r><w:r><w:t xml:space="preserve">ea, Chris, </w:t></w:r><w:r><w:t xml:space="preserve">Robert, Sergio &amp; </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>(Mab)</w:t></w:r><w:r><w:t>L</w:t></w:r><w:r><w:t>ead</w:t></w:r><w:r><w:t xml:space="preserve">er board separate page not working as </w:t></w:r><w:r><w:t>cannot send data to another domain</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>(Chris)Can’t randomise without losing submit, next and previous buttons</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">(Sergio)Presentation to Jawad, he is happy with progress </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>STARL not marked yet</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Next and Back buttons working (Andreea)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Possibly add leader board as part of questions page</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Need to start commenting out code</w:t></w:r><w:r><w:t xml:space="preserve"> (all of it including authors, versions etc)</w:t></w:r></w:p><w:p/>
'@

[void]$insertRange.InsertXML($newContentXml)
